$wb = $excel.ActiveWorkbook

# Work on the "protocoltestcasedetails" sheet that holds the test-case
# execution matrix (column D = execute Y/N flag).
$ws = $wb.Worksheets.Item("protocoltestcasedetails")

# Update the "execute" flag (Y/N) for several test cases:
#  - testcase18_parquet_dbtable_match_likeobject  -> no longer executed
#  - testcase19_oracle_mysql_match_manual          -> no longer executed
#  - testcase21_mysql_csv_match                    -> no longer executed
#  - testcase28_manual_sql_notifications           -> now executed
#  - testcase29_manual_sql_fullname                -> now executed
$ws.Range("D19").Value = "N"
$ws.Range("D20").Value = "N"
$ws.Range("D22").Value = "N"
$ws.Range("D29").Value = "Y"
$ws.Range("D30").Value = "Y"

# Update the sheet's stored view state: scroll back to the top of the
# sheet and move the active selection to C31.
$ws.Activate()
$ws.Range("C31").Select() | Out-Null
